# Natmi following Dr Hou advice
# Update LR-pair results for Ncam1-Gfra1: adds a new "ECs" sending-cluster
# category (3rd cluster alongside FAPs and sCs) and recomputes all
# statistics for rows 2-7 accordingly (row count grows from 4 to 6 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ncam1"
$ws.Cells.Item(2, 3).Value = "Gfra1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.9404873333333333
$ws.Cells.Item(2, 8).Value = 2.821462
$ws.Cells.Item(2, 9).Value = 0.02000383747045655
$ws.Cells.Item(2, 10).Value = 0.02000383747045654
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 8.922308333333334
$ws.Cells.Item(2, 14).Value = 26.766925
$ws.Cells.Item(2, 15).Value = 0.6138261687668722
$ws.Cells.Item(2, 16).Value = 0.6138261687668722
$ws.Cells.Item(2, 17).Value = 8.391317971594445
$ws.Cells.Item(2, 18).Value = 75.52186174435
$ws.Cells.Item(2, 19).Value = 0.01227887891512554
$ws.Cells.Item(2, 20).Value = 0.01227887891512554

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ncam1"
$ws.Cells.Item(3, 3).Value = "Gfra1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.9404873333333333
$ws.Cells.Item(3, 8).Value = 2.821462
$ws.Cells.Item(3, 9).Value = 0.02000383747045655
$ws.Cells.Item(3, 10).Value = 0.02000383747045654
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.613253666666666
$ws.Cells.Item(3, 14).Value = 16.839761
$ws.Cells.Item(3, 15).Value = 0.3861738312331279
$ws.Cells.Item(3, 16).Value = 0.3861738312331279
$ws.Cells.Item(3, 17).Value = 5.279193972286889
$ws.Cells.Item(3, 18).Value = 47.512745750582
$ws.Cells.Item(3, 19).Value = 0.007724958555331006
$ws.Cells.Item(3, 20).Value = 0.007724958555331004

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ncam1"
$ws.Cells.Item(4, 3).Value = "Gfra1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.392600333333333
$ws.Cells.Item(4, 8).Value = 4.177801000000001
$ws.Cells.Item(4, 9).Value = 0.0296201232509638
$ws.Cells.Item(4, 10).Value = 0.0296201232509638
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 8.922308333333334
$ws.Cells.Item(4, 14).Value = 26.766925
$ws.Cells.Item(4, 15).Value = 0.6138261687668722
$ws.Cells.Item(4, 16).Value = 0.6138261687668722
$ws.Cells.Item(4, 17).Value = 12.42520955910278
$ws.Cells.Item(4, 18).Value = 111.826886031925
$ws.Cells.Item(4, 19).Value = 0.01818160677354166
$ws.Cells.Item(4, 20).Value = 0.01818160677354166

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ncam1"
$ws.Cells.Item(5, 3).Value = "Gfra1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.392600333333333
$ws.Cells.Item(5, 8).Value = 4.177801000000001
$ws.Cells.Item(5, 9).Value = 0.0296201232509638
$ws.Cells.Item(5, 10).Value = 0.0296201232509638
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 5.613253666666666
$ws.Cells.Item(5, 14).Value = 16.839761
$ws.Cells.Item(5, 15).Value = 0.3861738312331279
$ws.Cells.Item(5, 16).Value = 0.3861738312331279
$ws.Cells.Item(5, 17).Value = 7.817018927284556
$ws.Cells.Item(5, 18).Value = 70.35317034556101
$ws.Cells.Item(5, 19).Value = 0.01143851647742214
$ws.Cells.Item(5, 20).Value = 0.01143851647742214

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Ncam1"
$ws.Cells.Item(6, 3).Value = "Gfra1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 44.682258
$ws.Cells.Item(6, 8).Value = 134.046774
$ws.Cells.Item(6, 9).Value = 0.9503760392785797
$ws.Cells.Item(6, 10).Value = 0.9503760392785796
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 8.922308333333334
$ws.Cells.Item(6, 14).Value = 26.766925
$ws.Cells.Item(6, 15).Value = 0.6138261687668722
$ws.Cells.Item(6, 16).Value = 0.6138261687668722
$ws.Cells.Item(6, 17).Value = 398.66888290555
$ws.Cells.Item(6, 18).Value = 3588.01994614995
$ws.Cells.Item(6, 19).Value = 0.583365683078205
$ws.Cells.Item(6, 20).Value = 0.583365683078205

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Ncam1"
$ws.Cells.Item(7, 3).Value = "Gfra1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 44.682258
$ws.Cells.Item(7, 8).Value = 134.046774
$ws.Cells.Item(7, 9).Value = 0.9503760392785797
$ws.Cells.Item(7, 10).Value = 0.9503760392785796
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 5.613253666666666
$ws.Cells.Item(7, 14).Value = 16.839761
$ws.Cells.Item(7, 15).Value = 0.3861738312331279
$ws.Cells.Item(7, 16).Value = 0.3861738312331279
$ws.Cells.Item(7, 17).Value = 250.812848553446
$ws.Cells.Item(7, 18).Value = 2257.315636981014
$ws.Cells.Item(7, 19).Value = 0.3670103562003748
$ws.Cells.Item(7, 20).Value = 0.3670103562003747
